$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-10-28 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-10-29 Tuesday", 2)

# Update each table cell with the new arithmetic expression
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "67-8="
$t.Cell(1, 2).Range.Text = "59-36="
$t.Cell(1, 3).Range.Text = "86-82="
$t.Cell(1, 4).Range.Text = "60+31="
$t.Cell(1, 5).Range.Text = "80-52="
$t.Cell(2, 1).Range.Text = "47+40="
$t.Cell(2, 2).Range.Text = "48-20="
$t.Cell(2, 3).Range.Text = "12+66="
$t.Cell(2, 4).Range.Text = "85-49="
$t.Cell(2, 5).Range.Text = "56-7="
$t.Cell(3, 1).Range.Text = "68+21="
$t.Cell(3, 2).Range.Text = "19+32="
$t.Cell(3, 3).Range.Text = "95-87="
$t.Cell(3, 4).Range.Text = "67-27="
$t.Cell(3, 5).Range.Text = "45-23="
$t.Cell(4, 1).Range.Text = "38+46="
$t.Cell(4, 2).Range.Text = "99-71="
$t.Cell(4, 3).Range.Text = "22+23="
$t.Cell(4, 4).Range.Text = "63-14="
$t.Cell(4, 5).Range.Text = "6+85="
$t.Cell(5, 1).Range.Text = "30+19="
$t.Cell(5, 2).Range.Text = "54-40="
$t.Cell(5, 3).Range.Text = "31+49="
$t.Cell(5, 4).Range.Text = "65-12="
$t.Cell(5, 5).Range.Text = "39-27="
$t.Cell(6, 1).Range.Text = "98-88="
$t.Cell(6, 2).Range.Text = "68-21="
$t.Cell(6, 3).Range.Text = "60-13="
$t.Cell(6, 4).Range.Text = "53+32="
$t.Cell(6, 5).Range.Text = "71-54="
$t.Cell(7, 1).Range.Text = "12+33="
$t.Cell(7, 2).Range.Text = "2+95="
$t.Cell(7, 3).Range.Text = "29-21="
$t.Cell(7, 4).Range.Text = "58-8="
$t.Cell(7, 5).Range.Text = "98-50="
$t.Cell(8, 1).Range.Text = "6+77="
$t.Cell(8, 2).Range.Text = "14+51="
$t.Cell(8, 3).Range.Text = "53-13="
$t.Cell(8, 4).Range.Text = "72+16="
$t.Cell(8, 5).Range.Text = "62+30="
$t.Cell(9, 1).Range.Text = "16+51="
$t.Cell(9, 2).Range.Text = "20-15="
$t.Cell(9, 3).Range.Text = "16+19="
$t.Cell(9, 4).Range.Text = "36+60="
$t.Cell(9, 5).Range.Text = "71+7="
$t.Cell(10, 1).Range.Text = "58-1="
$t.Cell(10, 2).Range.Text = "64+4="
$t.Cell(10, 3).Range.Text = "11+18="
$t.Cell(10, 4).Range.Text = "72-70="
$t.Cell(10, 5).Range.Text = "14+12="
$t.Cell(11, 1).Range.Text = "80-7="
$t.Cell(11, 2).Range.Text = "40+38="
$t.Cell(11, 3).Range.Text = "81-17="
$t.Cell(11, 4).Range.Text = "96-30="
$t.Cell(11, 5).Range.Text = "99-26="
$t.Cell(12, 1).Range.Text = "87-12="
$t.Cell(12, 2).Range.Text = "51-44="
$t.Cell(12, 3).Range.Text = "60-51="
$t.Cell(12, 4).Range.Text = "90-3="
$t.Cell(12, 5).Range.Text = "37+36="
$t.Cell(13, 1).Range.Text = "65-49="
$t.Cell(13, 2).Range.Text = "99-79="
$t.Cell(13, 3).Range.Text = "26+18="
$t.Cell(13, 4).Range.Text = "56+9="
$t.Cell(13, 5).Range.Text = "55+0="
$t.Cell(14, 1).Range.Text = "66+31="
$t.Cell(14, 2).Range.Text = "36+57="
$t.Cell(14, 3).Range.Text = "19+44="
$t.Cell(14, 4).Range.Text = "17+52="
$t.Cell(14, 5).Range.Text = "48+45="
$t.Cell(15, 1).Range.Text = "61-3="
$t.Cell(15, 2).Range.Text = "9+86="
$t.Cell(15, 3).Range.Text = "38-3="
$t.Cell(15, 4).Range.Text = "13+24="
$t.Cell(15, 5).Range.Text = "68-17="
$t.Cell(16, 1).Range.Text = "7+20="
$t.Cell(16, 2).Range.Text = "78-43="
$t.Cell(16, 3).Range.Text = "85-22="
$t.Cell(16, 4).Range.Text = "63-35="
$t.Cell(16, 5).Range.Text = "61+7="
$t.Cell(17, 1).Range.Text = "4+72="
$t.Cell(17, 2).Range.Text = "87-19="
$t.Cell(17, 3).Range.Text = "30-0="
$t.Cell(17, 4).Range.Text = "59+24="
$t.Cell(17, 5).Range.Text = "78-46="
$t.Cell(18, 1).Range.Text = "95+0="
$t.Cell(18, 2).Range.Text = "79-23="
$t.Cell(18, 3).Range.Text = "32+48="
$t.Cell(18, 4).Range.Text = "36+43="
$t.Cell(18, 5).Range.Text = "82-16="
$t.Cell(19, 1).Range.Text = "29+67="
$t.Cell(19, 2).Range.Text = "69+1="
$t.Cell(19, 3).Range.Text = "22+62="
$t.Cell(19, 4).Range.Text = "72-61="
$t.Cell(19, 5).Range.Text = "37+35="
$t.Cell(20, 1).Range.Text = "39-18="
$t.Cell(20, 2).Range.Text = "88-76="
$t.Cell(20, 3).Range.Text = "5+80="
$t.Cell(20, 4).Range.Text = "70+24="
$t.Cell(20, 5).Range.Text = "11+70="
